$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# A1 stays "BPMN_File_Name" (unchanged). B1:D1 get new metric names,
# and E1:L1 are new columns with more metric names.
$headers = @("nTask", "nSendTask", "nUserTask", "nManualTask", "nBusinessRuleTask", "nServiceTask", "nScriptTask", "nCallActivity", "nSubProcess", "nTransaction", "nAdHocSubProcess")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Copy the header style (bold + border + centered) from A1 onto the newly
# added header cells E1:L1 so the whole header row is formatted uniformly.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1:L1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data row (row 2) ---
$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = 1
for ($col = 3; $col -le 12; $col++) {
    $ws.Cells.Item(2, $col).Value = 0
}
